$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.005.71"
$ws.Range("E2").Value = "  +0.48%  "
$ws.Range("D3").Value = "3.777.90"
$ws.Range("E3").Value = "  -0.64%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'598.63"
$ws.Range("E5").Value = "  +0.27%  "
$ws.Range("D6").Value = "'163.31"
$ws.Range("E6").Value = "  -2.33%  "
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("D8").Value = "'0.514"
$ws.Range("E8").Value = "  -1.02%  "
$ws.Range("D9").Value = "'0.157"
$ws.Range("E9").Value = "  -1.71%  "
$ws.Range("D10").Value = "'0.447"
$ws.Range("E10").Value = "  -0.65%  "
$ws.Range("D11").Value = "'6.54"
$ws.Range("E11").Value = "  +3.83%  "
$ws.Range("D12").Value = "'0.0000246"
$ws.Range("E12").Value = "  -2.76%  "
$ws.Range("D13").Value = "'35.38"
$ws.Range("E13").Value = "  -1.44%  "
$ws.Range("D14").Value = "4.409.58"
$ws.Range("E14").Value = "  -0.74%  "
$ws.Range("D15").Value = "3.774.26"
$ws.Range("E15").Value = "  -0.79%  "
$ws.Range("D16").Value = "67.974.65"
$ws.Range("E16").Value = "  +0.39%  "
$ws.Range("D17").Value = "'18.23"
$ws.Range("E17").Value = "  -2.17%  "
$ws.Range("E18").Value = "  +2.04%  "
$ws.Range("D19").Value = "'7.00"
$ws.Range("E19").Value = "  -1.25%  "
$ws.Range("D20").Value = "'458.73"
$ws.Range("E20").Value = "  -0.54%  "
$ws.Range("D21").Value = "'9.61"
$ws.Range("E21").Value = "  -3.32%  "
$ws.Range("D22").Value = "'0.696"
$ws.Range("E22").Value = "  -0.91%  "
$ws.Range("D23").Value = "'82.77"
$ws.Range("E23").Value = "  -0.78%  "
$ws.Range("D24").Value = "'0.0000143"
$ws.Range("E24").Value = "  -5.89%  "
$ws.Range("E25").Value = "  -1.50%  "
$ws.Range("D26").Value = "'2.08"
$ws.Range("E26").Value = "  -0.86%  "
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("D28").Value = "'9.92"
$ws.Range("E28").Value = "  -0.88%  "
$ws.Range("D29").Value = "3.929.54"
$ws.Range("E29").Value = "  -0.38%  "
$ws.Range("E30").Value = "  -0.86%  "
$ws.Range("D31").Value = "'7.22"
$ws.Range("E31").Value = "  -0.93%  "
$ws.Range("D32").Value = "'2.55"
$ws.Range("E32").Value = "  -7.85%  "
$ws.Range("D33").Value = "'29.00"
$ws.Range("E33").Value = "  -2.07%  "
$ws.Range("E34").Value = "  +0.12%  "
$ws.Range("D35").Value = "'8.93"
$ws.Range("E35").Value = "  -1.46%  "
$ws.Range("D36").Value = "'0.0991"
$ws.Range("E36").Value = "  -0.84%  "
$ws.Range("D37").Value = "'0.141"
$ws.Range("E37").Value = "  +2.43%  "
$ws.Range("D38").Value = "'5.78"
$ws.Range("E38").Value = "  +0.04%  "
$ws.Range("D39").Value = "'3.21"
$ws.Range("E39").Value = "  -5.03%  "
$ws.Range("D40").Value = "'0.980"
$ws.Range("E40").Value = "  -1.92%  "
$ws.Range("D41").Value = "'0.999"
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("B43").Value = "Arweave"
$ws.Range("C43").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D43").Value = "'43.45"
$ws.Range("E43").Value = "  +0.98%  "
$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").Value = "'47.23"
$ws.Range("E44").Value = "  -1.92%  "
$ws.Range("B45").Value = "Monero"
$ws.Range("C45").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D45").Value = "'152.74"
$ws.Range("E45").Value = "  +3.66%  "
$ws.Range("B46").Value = "TheGraph"
$ws.Range("C46").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D46").Value = "'0.296"
$ws.Range("E46").Value = "  -1.55%  "
$ws.Range("D47").Value = "'8.29"
$ws.Range("E47").Value = "  -0.50%  "
$ws.Range("D48").Value = "'1.37"
$ws.Range("E48").Value = "  +0.26%  "
$ws.Range("E49").Value = "  +0.14%  "
$ws.Range("D50").Value = "'387.17"
$ws.Range("E50").Value = "  -1.84%  "
$ws.Range("D51").Value = "'26.66"
$ws.Range("E51").Value = "  -1.96%  "
